$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.619.35"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.698.43"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "672.98"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.38"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.497"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("E9").Value = "  +0.65%  "
$ws.Range("E10").Value = "  +0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.444"
$ws.Range("E11").Value = "  +2.35%  "
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("E13").Value = "  +1.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.708.93"
$ws.Range("E14").Value = "  +1.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.645.46"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("E16").Value = "  +1.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.16"
$ws.Range("E17").Value = "  +2.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.50"
$ws.Range("E18").Value = "  +2.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "473.94"
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.654"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.47"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.846.99"
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.03"
$ws.Range("E26").Value = "  +0.71%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.12"
$ws.Range("E27").Value = "  -0.72%  "
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  +0.42%  "
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.60"
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("E32").Value = "  +0.22%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.88"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.686.59"
$ws.Range("E35").Value = "  +1.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.50"
$ws.Range("E36").Value = "  +4.26%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "172.53"
$ws.Range("E42").Value = "  +3.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.940"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("E45").Value = "  +2.36%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.05"
$ws.Range("E46").Value = "  +1.52%  "
$ws.Range("B47").Value = "FLOKI"
$ws.Range("C47").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000278"
$ws.Range("E47").Value = "  -0.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.11"
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  -0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.89"
$ws.Range("E50").Value = "  +1.74%  "
$ws.Range("E51").Value = "  +0.77%  "
